$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve a copy of the "Hipervinculo" (hyperlink) cell look-and-feel in a
# scratch cell before we touch anything, so we can re-apply it later without
# Excel fabricating brand-new style records for every Hyperlinks.Add() call.
$ws.Range("N2").Copy()
$ws.Range("Z100").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- Remove the two existing hyperlinks (N2 mailto, V2 course url) ---
$existing = @($ws.Hyperlinks)
for ($i = $existing.Length - 1; $i -ge 0; $i--) {
    $existing[$i].Delete()
}

# N2 no longer is a hyperlink -> drop its "Hipervinculo" look, reuse the plain
# look already used by the rest of the row (M2).
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the row 2 values ---
$ws.Range("A2").Value = 150
$ws.Range("E2").Value = "FILO_7014"
$ws.Range("F2").Value = "Bioética"
$ws.Range("G2").Value = "Vallejo Delgado Merci Lorena"
$ws.Range("H2").Value = "mlvallejo11@utpl.edu.ec"
$ws.Range("I2").Value = 1102094883
$ws.Range("K2").Value = "Ciencias de la Salud"
$ws.Range("L2").Value = "Maestría en Gerencia de Instituciones de Salud"
$ws.Range("M2").Value = "Gloria Alexandra Carrión Figueroa"
$ws.Range("N2").Value = "gacarrionx@utpl.edu.ec"
$ws.Range("Q2").Value = "Total 144: ACD_32 APE_24 AA_88"
$ws.Range("S2").Value = 4
$ws.Range("T2").Value = "Trayectoria profesional"
$ws.Range("U2").Value = "SI"
$ws.Range("X2").Value = "FILO_7014_NLG_META"
$ws.Range("Y2").Value = ""
$ws.Range("AA2").Value = "FILO_7014"
$ws.Range("AB2").Value = 45835
$ws.Range("AC2").Value = 45835

# --- Re-create the hyperlinks ---
# H2: new hyperlink, whose "display" attribute equals the address, while the
# visible cell text stays the docent's e-mail (set above).
$addr = "http://mlvallejo11utpl.edu.ec/"
$ws.Hyperlinks.Add($ws.Range("H2"), $addr, [Type]::Missing, [Type]::Missing, $addr)
$ws.Range("H2").Value = "mlvallejo11@utpl.edu.ec"
$ws.Range("Z100").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# V2: new course url, no explicit display text (matches original behaviour).
$ws.Range("V2").Value = "https://utpl.instructure.com/courses/56724"
$ws.Hyperlinks.Add($ws.Range("V2"), "https://utpl.instructure.com/courses/56724")
$ws.Range("Z100").Copy()
$ws.Range("V2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# clean up scratch cell
$ws.Range("Z100").Clear()

# --- Update the view: scroll so column D is left-most, select V2 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("V2").Select() | Out-Null
